# Auto-generated: update Leve profit-calculation columns (H-N) across sheets
# per scheduled market-price refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 949.75
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 1749.5
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 1749.5
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -1975.5
$ws.Range("H6").Value = 2000574.9
$ws.Range("J6").Value = 4400
$ws.Range("L6").Value = 13200
$ws.Range("N6").Value = -13424
$ws.Range("H64").Value = 9749.25
$ws.Range("J64").Value = 9665.666999999999
$ws.Range("L64").Value = 9665.666999999999
$ws.Range("N64").Value = -10161.667
$ws.Range("H67").Value = 9749.25
$ws.Range("J67").Value = 9665.666999999999
$ws.Range("L67").Value = 9665.666999999999
$ws.Range("N67").Value = -11381.667
$ws.Range("H76").Value = 16560.625
$ws.Range("I76").Value = 16211.571
$ws.Range("K76").Value = 16211.571
$ws.Range("M76").Value = -15896.571
$ws.Range("H79").Value = 16560.625
$ws.Range("I79").Value = 16211.571
$ws.Range("K79").Value = 16211.571
$ws.Range("M79").Value = -15119.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H63").Value = 2836.2727
$ws.Range("I63").Value = 2537.5
$ws.Range("J63").Value = 3633
$ws.Range("K63").Value = 2537.5
$ws.Range("L63").Value = 3633
$ws.Range("M63").Value = -1851.5
$ws.Range("N63").Value = -5005
$ws.Range("H66").Value = 2836.2727
$ws.Range("I66").Value = 2537.5
$ws.Range("J66").Value = 3633
$ws.Range("K66").Value = 12687.5
$ws.Range("L66").Value = 18165
$ws.Range("M66").Value = -9255.5
$ws.Range("N66").Value = -25029
$ws.Range("H132").Value = 3098.2632
$ws.Range("I132").Value = 1616.4375
$ws.Range("K132").Value = 4849.3125
$ws.Range("M132").Value = -2319.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 42490.832
$ws.Range("I82").Value = 7472.5
$ws.Range("K82").Value = 7472.5
$ws.Range("M82").Value = -7089.5
$ws.Range("H85").Value = 42490.832
$ws.Range("I85").Value = 7472.5
$ws.Range("K85").Value = 7472.5
$ws.Range("M85").Value = -6146.5
$ws.Range("H86").Value = 4224.4
$ws.Range("I86").Value = 1913.25
$ws.Range("K86").Value = 1913.25
$ws.Range("M86").Value = -790.25
$ws.Range("H89").Value = 4224.4
$ws.Range("I89").Value = 1913.25
$ws.Range("K89").Value = 9566.25
$ws.Range("M89").Value = -3950.25
$ws.Range("H94").Value = 887.25
$ws.Range("J94").Value = 651.5
$ws.Range("L94").Value = 651.5
$ws.Range("N94").Value = -1553.5
$ws.Range("H97").Value = 2875.6667
$ws.Range("I97").Value = 2875.6667
$ws.Range("K97").Value = 2875.6667
$ws.Range("M97").Value = -1884.6667
$ws.Range("H107").Value = 1550.8
$ws.Range("I107").Value = 1395.3334
$ws.Range("K107").Value = 1395.3334
$ws.Range("M107").Value = 524.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1299.5238
$ws.Range("I16").Value = 1303.1666
$ws.Range("J16").Value = 1277.6666
$ws.Range("K16").Value = 1303.1666
$ws.Range("L16").Value = 1277.6666
$ws.Range("M16").Value = -1016.1666
$ws.Range("N16").Value = -1851.6666
$ws.Range("H76").Value = 6025
$ws.Range("I76").Value = 6025
$ws.Range("K76").Value = 6025
$ws.Range("M76").Value = -5710
$ws.Range("H79").Value = 6025
$ws.Range("I79").Value = 6025
$ws.Range("K79").Value = 6025
$ws.Range("M79").Value = -4933
$ws.Range("H113").Value = 1299.5238
$ws.Range("I113").Value = 1303.1666
$ws.Range("J113").Value = 1277.6666
$ws.Range("K113").Value = 1303.1666
$ws.Range("L113").Value = 1277.6666
$ws.Range("M113").Value = 866.8334
$ws.Range("N113").Value = -5617.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 60304.5
$ws.Range("I2").Value = 25.5
$ws.Range("J2").Value = 75374.25
$ws.Range("K2").Value = 153
$ws.Range("L2").Value = 452245.5
$ws.Range("M2").Value = -40
$ws.Range("N2").Value = -452471.5
$ws.Range("H39").Value = 3108.0833
$ws.Range("I39").Value = 1823.75
$ws.Range("J39").Value = 3750.25
$ws.Range("K39").Value = 5471.25
$ws.Range("L39").Value = 11250.75
$ws.Range("M39").Value = -5177.25
$ws.Range("N39").Value = -11838.75
$ws.Range("H56").Value = 5998.5557
$ws.Range("I56").Value = 5998.5557
$ws.Range("K56").Value = 5998.5557
$ws.Range("M56").Value = -5468.5557
$ws.Range("H75").Value = 4468.9
$ws.Range("I75").Value = 458
$ws.Range("J75").Value = 6187.857
$ws.Range("K75").Value = 1374
$ws.Range("L75").Value = 18563.571
$ws.Range("M75").Value = -376
$ws.Range("N75").Value = -20559.571
$ws.Range("H78").Value = 4468.9
$ws.Range("I78").Value = 458
$ws.Range("J78").Value = 6187.857
$ws.Range("K78").Value = 4122
$ws.Range("L78").Value = 55690.713
$ws.Range("M78").Value = 870
$ws.Range("N78").Value = -65674.713
$ws.Range("H92").Value = 3033.8333
$ws.Range("J92").Value = 7501.5
$ws.Range("L92").Value = 22504.5
$ws.Range("N92").Value = -25000.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 31394.5
$ws.Range("J55").Value = 31394.5
$ws.Range("L55").Value = 31394.5
$ws.Range("N55").Value = -32048.5
$ws.Range("H58").Value = 52543.332
$ws.Range("I58").Value = 52543.332
$ws.Range("K58").Value = 52543.332
$ws.Range("M58").Value = -52266.332
$ws.Range("H92").Value = 14537.75
$ws.Range("J92").Value = 14537.75
$ws.Range("L92").Value = 14537.75
$ws.Range("N92").Value = -18281.75
$ws.Range("H102").Value = 3892.0588
$ws.Range("I102").Value = 1903.3846
$ws.Range("K102").Value = 1903.3846
$ws.Range("M102").Value = -281.3846000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5232.385
$ws.Range("I40").Value = 3651.2
$ws.Range("J40").Value = 10503
$ws.Range("K40").Value = 3651.2
$ws.Range("L40").Value = 10503
$ws.Range("M40").Value = -3515.2
$ws.Range("N40").Value = -10775

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2249.5
$ws.Range("I62").Value = 1499
$ws.Range("K62").Value = 1499
$ws.Range("M62").Value = -875
$ws.Range("H65").Value = 2249.5
$ws.Range("I65").Value = 1499
$ws.Range("K65").Value = 7495
$ws.Range("M65").Value = -4375
$ws.Range("H126").Value = 1948.6428
$ws.Range("I126").Value = 1114.7
$ws.Range("K126").Value = 3344.1
$ws.Range("M126").Value = -874.1000000000004
$ws.Range("H136").Value = 4085.1155
$ws.Range("I136").Value = 2012.3889
$ws.Range("K136").Value = 6037.1667
$ws.Range("M136").Value = -3487.1667

Write-Host "Updated Leve profit columns across all sheets."